# Add "NON-FUNCTIONAL REQUIREMENTS" section to the requirements sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Section header - bold, otherwise same font as the rest of the sheet.
$ws.Range("A28").Value = "NON-FUNCTIONAL REQUIREMENTS"
$ws.Range("A28").Font.Name = "Arial"
$ws.Range("A28").Font.Size = 10
$ws.Range("A28").Font.Bold = $true

# Individual non-functional requirement bullets.
$ws.Range("A29").Value = "Should have an easy to navigate and attractive user interface"

$ws.Range("A33").Value = "Should be responsive and work on phones as well"

$ws.Range("A37").Value = "Should work in all commonly used browsers"
$ws.Range("A37").Font.Name = "Arial"
$ws.Range("A37").Font.Size = 10
$ws.Range("A37").Font.Bold = $false

$ws.Range("A41").Value = "Performance should be good enough so that it does not take too long to load the recipes"
$ws.Range("A41").Font.Name = "Arial"
$ws.Range("A41").Font.Size = 10
$ws.Range("A41").Font.Bold = $false

$ws.Range("A46").Value = "Storage and processing of user info must be secure"
$ws.Range("A46").Font.Name = "Arial"
$ws.Range("A46").Font.Size = 10
$ws.Range("A46").Font.Bold = $false

# Leave the selection on the last cell typed, matching the author's final
# cursor position when they saved the workbook.
$ws.Range("A46").Select() | Out-Null
